$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Link/Price/Volume cells are text-formatted in the source data; Price
# values that look numeric (e.g. "21.33") need an explicit Text number format
# so Excel does not silently coerce them to the Number type on assignment.

$ws.Range("D2").Value = '71.038.11'
$ws.Range("E2").Value = '  -2.31%  '

$ws.Range("D3").Value = '3.870.05'
$ws.Range("E3").Value = '  -2.33%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.12'
$ws.Range("E5").Value = '  +1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.83'
$ws.Range("E6").Value = '  +5.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.672'
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.753'
$ws.Range("E9").Value = '  +0.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.176'
$ws.Range("E10").Value = '  +4.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.47'
$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000321'
$ws.Range("E12").Value = '  +0.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.44'
$ws.Range("E13").Value = '  +5.52%  '

$ws.Range("D14").Value = '4.480.40'
$ws.Range("E14").Value = '  -2.73%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.33'
$ws.Range("E15").Value = '  +4.98%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.874.88'
$ws.Range("E16").Value = '  -2.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.84'
$ws.Range("E17").Value = '  -1.38%  '

$ws.Range("E18").Value = '  -4.24%  '

$ws.Range("E19").Value = '  -2.09%  '

$ws.Range("D20").Value = '70.798.31'
$ws.Range("E20").Value = '  -2.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '437.08'
$ws.Range("E21").Value = '  +0.77%  '

$ws.Range("E22").Value = '  +0.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '94.45'
$ws.Range("E23").Value = '  -1.41%  '

$ws.Range("E24").Value = '  -4.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.88'
$ws.Range("E25").Value = '  -3.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.37'
$ws.Range("E26").Value = '  +2.33%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.08'
$ws.Range("E27").Value = '  -7.49%  '

$ws.Range("E28").Value = '  -0.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.40'
$ws.Range("E29").Value = '  -1.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.19'
$ws.Range("E30").Value = '  -3.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.21'
$ws.Range("E31").Value = '  +4.34%  '

$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '48.25'
$ws.Range("E33").Value = '  -0.47%  '

$ws.Range("E34").Value = '  -4.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '69.96'
$ws.Range("E35").Value = '  +0.56%  '

$ws.Range("D36").Value = '0.0₃0985'
$ws.Range("E36").Value = '  +13.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '630.75'
$ws.Range("E37").Value = '  -7.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.428'
$ws.Range("E38").Value = '  -1.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.148'
$ws.Range("E39").Value = '  +0.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.30'
$ws.Range("E41").Value = '  -2.43%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.28'
$ws.Range("E43").Value = '  +26.37%  '

$ws.Range("E44").Value = '  -3.30%  '

$ws.Range("E45").Value = '  -6.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.70'
$ws.Range("E46").Value = '  +1.86%  '

$ws.Range("E47").Value = '  -3.37%  '

$ws.Range("E48").Value = '  -1.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("E49").Value = '  -15.17%  '

$ws.Range("D50").Value = '2.844.90'
$ws.Range("E50").Value = '  +2.03%  '

$ws.Range("E51").Value = '  +1.58%  '
